# Updated capital structure database
# Applies refreshed metric values to rows 2 and 3 (Pakistan / Engineering-Construction
# industry-average row and Gammon Pakistan Limited row) and removes the now-obsolete
# historical_growth_net_income_last_5_years (column E) figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2, 3) {
    $ws.Range("D$row").Value = -0.0806
    $ws.Range("E$row").ClearContents()
    $ws.Range("G$row").Value = -0.01219512195121951
    $ws.Range("H$row").Value = -0.01219512195121951
    $ws.Range("I$row").Value = -0.03089430894308943
    $ws.Range("J$row").Value = -0.01800315844885068
    $ws.Range("K$row").Value = 0.162
    $ws.Range("L$row").Value = 0.1317073170731707
    $ws.Range("U$row").Value = 0.005
    $ws.Range("V$row").Value = 0.002173913043478261
    $ws.Range("W$row").Value = 0.03454157782515991
    $ws.Range("X$row").Value = 0.09379303256626192
    $ws.Range("Y$row").Value = -0.05925145474110202
    $ws.Range("Z$row").Value = 0.2514309076042518
    $ws.Range("AA$row").Value = -0.00452655046853768
    $ws.Range("AB$row").Value = 0.09379303256626192
    $ws.Range("AC$row").Value = -0.09831958303479961
    $ws.Range("AD$row").Value = 0
    $ws.Range("AF$row").Value = 0
    $ws.Range("AG$row").Value = -0.005
    $ws.Range("AH$row").Value = 0
    $ws.Range("AI$row").Value = 0
    $ws.Range("AJ$row").Value = -0.002178649237472767
    $ws.Range("AK$row").Value = -0.001055966209081309
    $ws.Range("AN$row").Value = 0
    $ws.Range("AP$row").Value = 0.3571428571428572
}
